$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 165, shifting existing rows 165:236 down to 166:237
$ws.Rows(165).Insert()

# Populate the newly inserted row 165 with the new weekly record
$ws.Range("A165").Value = 7
$ws.Range("B165").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C165").Value = "Ñuble"
$ws.Range("D165").Value = 44704
$ws.Range("E165").Value = 16
$ws.Range("F165").Value = 100112003
$ws.Range("G165").Value = "Ajo"
$ws.Range("H165").Value = "Chino"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 60
$ws.Range("K165").Value = 19000
$ws.Range("L165").Value = 20000
$ws.Range("M165").Value = 19500
$ws.Range("N165").Value = '$/caja 10 kilos'
$ws.Range("O165").Value = "China"
$ws.Range("P165").Value = 1950
$ws.Range("Q165").Value = 10
$ws.Range("R165").Value = "Hortaliza"
